$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`r`nWITH DISTINCT f, parent`r`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`r`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`r`n MATCH (samp:sample)-->(c) `r`n WHERE f.file_format IN [`"doc`"]  `r`nWITH DISTINCT f, parent, c, demo, diag, s`r`nRETURN coalesce(f.file_name, '') AS ``File Name``, `r`n        coalesce(labels(parent)[0], '') AS ``Association``,`r`n        coalesce(f.file_description, '') AS ``Description``,`r`n        coalesce(f.file_format, '') AS ``Format``,`r`n        coalesce(f.file_size, '') AS ``Size``,`r`n        coalesce(c.case_id, '') AS ``Case ID``, `r`n        coalesce(diag.disease_term,'') AS Diagnosis , `r`n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Select()
